$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Duration_Seconds (column C) values per the new protocol timings
$ws.Range("C3").Value = 10
$ws.Range("C4").Value = 0.5
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 0.5
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 0.5
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 10
$ws.Range("C12").Value = 0.1
$ws.Range("C13").Value = 0.1
$ws.Range("C14").Value = 0.1
$ws.Range("C15").Value = 0.1
$ws.Range("C16").Value = 0.1
$ws.Range("C17").Value = 0.1
$ws.Range("C18").Value = 0.1
$ws.Range("C19").Value = 0.1
$ws.Range("C20").Value = 0.1
$ws.Range("C21").Value = 0.1
$ws.Range("C22").Value = 0.1
$ws.Range("C23").Value = 0.1
$ws.Range("C24").Value = 0.1
$ws.Range("C25").Value = 0.1
$ws.Range("C26").Value = 0.1
$ws.Range("C27").Value = 0.1
$ws.Range("C28").Value = 0.1
$ws.Range("C29").Value = 0.1
$ws.Range("C30").Value = 0.1
$ws.Range("C31").Value = 0.1
$ws.Range("C32").Value = 0.1
$ws.Range("C33").Value = 10
$ws.Range("C34").Value = 0.1
$ws.Range("C35").Value = 10
$ws.Range("C36").Value = 0.1
$ws.Range("C37").Value = 10
$ws.Range("C38").Value = 0.1
$ws.Range("C39").Value = 10
$ws.Range("C40").Value = 0.1
$ws.Range("C41").Value = 10
$ws.Range("C42").Value = 1

# Move the active selection (cursor was left on C43 after editing the sheet)
$ws.Range("C43").Select()
